# Weekly update: insert a new week of "Cebolla" price records
# (date serial 44753) for "Femacal de La Calera", pushing the existing
# rows 1295-1350 down by three rows (new dimension A1:R1353).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows right before the current row 1295; everything
# that used to live at row 1295 onward shifts down to row 1298 onward.
$ws.Rows("1295:1297").Insert()

# --- New row 1295 --------------------------------------------------
$ws.Cells.Item(1295, 1).Value  = 3
$ws.Cells.Item(1295, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(1295, 3).Value  = "Coquimbo"
$ws.Cells.Item(1295, 4).Value  = 44753
$ws.Cells.Item(1295, 5).Value  = 5
$ws.Cells.Item(1295, 6).Value  = 100112004
$ws.Cells.Item(1295, 7).Value  = "Cebolla"
$ws.Cells.Item(1295, 8).Value  = "Morada(o)"
$ws.Cells.Item(1295, 9).Value  = "Primera"
$ws.Cells.Item(1295, 10).Value = 171
$ws.Cells.Item(1295, 11).Value = 9000
$ws.Cells.Item(1295, 12).Value = 12500
$ws.Cells.Item(1295, 13).Value = 11404
$ws.Cells.Item(1295, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(1295, 15).Value = "Perú"
$ws.Cells.Item(1295, 16).Value = 634
$ws.Cells.Item(1295, 17).Value = 18
$ws.Cells.Item(1295, 18).Value = "Hortaliza"

# --- New row 1296 --------------------------------------------------
$ws.Cells.Item(1296, 1).Value  = 3
$ws.Cells.Item(1296, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(1296, 3).Value  = "Coquimbo"
$ws.Cells.Item(1296, 4).Value  = 44753
$ws.Cells.Item(1296, 5).Value  = 5
$ws.Cells.Item(1296, 6).Value  = 100112004
$ws.Cells.Item(1296, 7).Value  = "Cebolla"
$ws.Cells.Item(1296, 8).Value  = "Sin especificar"
$ws.Cells.Item(1296, 9).Value  = "1a (guarda)"
$ws.Cells.Item(1296, 10).Value = 175
$ws.Cells.Item(1296, 11).Value = 6000
$ws.Cells.Item(1296, 12).Value = 6300
$ws.Cells.Item(1296, 13).Value = 6137
$ws.Cells.Item(1296, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(1296, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1296, 16).Value = 341
$ws.Cells.Item(1296, 17).Value = 18
$ws.Cells.Item(1296, 18).Value = "Hortaliza"

# --- New row 1297 --------------------------------------------------
$ws.Cells.Item(1297, 1).Value  = 3
$ws.Cells.Item(1297, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(1297, 3).Value  = "Coquimbo"
$ws.Cells.Item(1297, 4).Value  = 44753
$ws.Cells.Item(1297, 5).Value  = 5
$ws.Cells.Item(1297, 6).Value  = 100112004
$ws.Cells.Item(1297, 7).Value  = "Cebolla"
$ws.Cells.Item(1297, 8).Value  = "Sin especificar"
$ws.Cells.Item(1297, 9).Value  = "2a (guarda)"
$ws.Cells.Item(1297, 10).Value = 50
$ws.Cells.Item(1297, 11).Value = 5000
$ws.Cells.Item(1297, 12).Value = 5000
$ws.Cells.Item(1297, 13).Value = 5000
$ws.Cells.Item(1297, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(1297, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1297, 16).Value = 278
$ws.Cells.Item(1297, 17).Value = 18
$ws.Cells.Item(1297, 18).Value = "Hortaliza"
